$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.656.65"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.038.29"
$ws.Range("E3").Value = "  +3.08%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.26"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.17"
$ws.Range("E6").Value = "  +1.94%  "

# Row 7
$ws.Range("E7").Value = "  +1.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +2.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.92"
$ws.Range("E10").Value = "  +1.99%  "

# Row 11
$ws.Range("E11").Value = "  -0.21%  "

# Row 12
$ws.Range("E12").Value = "  +1.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.517.08"
$ws.Range("E13").Value = "  +3.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.56"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.74"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.027.66"
$ws.Range("E16").Value = "  +2.69%  "

# Row 17
$ws.Range("E17").Value = "  -2.24%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.52"
$ws.Range("E18").Value = "  -11.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.669.19"
$ws.Range("E19").Value = "  +1.39%  "

# Row 20
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +1.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.94"
$ws.Range("E23").Value = "  +0.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.63"
$ws.Range("E24").Value = "  +0.82%  "

# Row 25
$ws.Range("E25").Value = "  -2.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.18"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.63"
$ws.Range("E27").Value = "  +7.88%  "

# Row 28
$ws.Range("E28").Value = "  +6.57%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.34"
$ws.Range("E29").Value = "  +2.89%  "

# Row 30
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.29"
$ws.Range("E32").Value = "  +2.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.19"
$ws.Range("E33").Value = "  +2.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.53"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35
$ws.Range("E35").Value = "  +0.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +4.72%  "

# Row 37
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +7.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.292"
$ws.Range("E39").Value = "  +12.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.13"
$ws.Range("E40").Value = "  +2.98%  "

# Row 41
$ws.Range("E41").Value = "  +3.40%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +2.77%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.34"
$ws.Range("E44").Value = "  +8.20%  "

# Row 45
$ws.Range("E45").Value = "  +6.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.78"
$ws.Range("E46").Value = "  +1.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").Value = "  +5.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +2.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.036.46"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.335.62"
$ws.Range("E50").Value = "  +2.90%  "

# Row 51
$ws.Range("E51").Value = "  +2.31%  "
